$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("19:19").Insert()

$ws.Range("A19").Value = 7
$ws.Range("B19").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C19").Value = "Ñuble"
$ws.Range("D19").Value = 44910
$ws.Range("E19").Value = 16
$ws.Range("F19").Value = 100112031
$ws.Range("G19").Value = "Poroto verde"
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 50
$ws.Range("K19").Value = 30000
$ws.Range("L19").Value = 30000
$ws.Range("M19").Value = 30000
$ws.Range("N19").Value = "$/saco 25 kilos"
$ws.Range("O19").Value = "Región del Maule"
$ws.Range("P19").Value = 1200
$ws.Range("Q19").Value = 25
$ws.Range("R19").Value = "Hortaliza"
